$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Range("A" + $row).Value = 'Raul Camargo_20251201_165414'
$ws.Range("B" + $row).Value = "'"
$ws.Range("B" + $row).Style = "Normal"
$ws.Range("C" + $row).Value = 'Raul Camargo'
$ws.Range("D" + $row).Value = 21
$ws.Range("E" + $row).Value = 'Male'
$ws.Range("F" + $row).Value = '2025-12-01 16:54:14'
$ws.Range("G" + $row).Value = '{
  "portion": 0.6,
  "diet": 0.2857142857142857,
  "salt": 0.2,
  "fat": 0.2,
  "natural": 0.8,
  "convenience": 0.2,
  "price": 0.2
}'
$ws.Range("H" + $row).Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("I" + $row).Value = '''0.721'
$ws.Range("I" + $row).Style = "Normal"
$ws.Range("J" + $row).Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("K" + $row).Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("L" + $row).Value = '''0.449'
$ws.Range("L" + $row).Style = "Normal"
$ws.Range("M" + $row).Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("N" + $row).Value = 'Nongshim Shin Ramyun'
$ws.Range("O" + $row).Value = '''0.418'
$ws.Range("O" + $row).Style = "Normal"
$ws.Range("P" + $row).Value = 'Sabor intenso, picante, umami, fideos gruesos, muy alto en sodio'
$ws.Range("Q" + $row).Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("R" + $row).Value = '''0.746'
$ws.Range("R" + $row).Style = "Normal"
$ws.Range("S" + $row).Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("T" + $row).Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("U" + $row).Value = '''0.670'
$ws.Range("U" + $row).Style = "Normal"
$ws.Range("V" + $row).Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("W" + $row).Value = 'Annie’s Shells & White Cheddar'
$ws.Range("X" + $row).Value = '''0.641'
$ws.Range("X" + $row).Style = "Normal"
$ws.Range("Y" + $row).Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("Z" + $row).Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("AA" + $row).Value = '''0.754'
$ws.Range("AA" + $row).Style = "Normal"
$ws.Range("AB" + $row).Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AC" + $row).Value = 'Kitchens of India Variety Pack'
$ws.Range("AD" + $row).Value = '''0.597'
$ws.Range("AD" + $row).Style = "Normal"
$ws.Range("AE" + $row).Value = 'Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad'
$ws.Range("AF" + $row).Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AG" + $row).Value = '''0.363'
$ws.Range("AG" + $row).Style = "Normal"
$ws.Range("AH" + $row).Value = 'Portátil, saludable, fácil, buena textura, sabor suave'

$ws.Rows.Item($row).AutoFit()

